$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price (D) column cells whose new values would
# otherwise be auto-interpreted as numbers by Excel, so they stay as plain
# text matching the source data, just like the other price cells that are
# already non-numeric-looking. (A contiguous range is used because this
# Excel engine only honors NumberFormat assignment on the first area of a
# comma-separated multi-area range.)
$ws.Range("D2:D51").NumberFormat = "@"

# Apply the updated cell values row by row

# Row 2
$ws.Range("D2").Value = '43.004.46'
$ws.Range("E2").Value = '  +0.00%  '

# Row 3
$ws.Range("D3").Value = '2.330.30'
$ws.Range("E3").Value = '  +1.08%  '

# Row 4
$ws.Range("E4").Value = '  +0.20%  '

# Row 5
$ws.Range("D5").Value = '303.30'
$ws.Range("E5").Value = '  -0.51%  '

# Row 6
$ws.Range("D6").Value = '96.04'
$ws.Range("E6").Value = '  -1.12%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("E8").Value = '  +0.21%  '

# Row 9
$ws.Range("D9").Value = '0.495'
$ws.Range("E9").Value = '  -1.09%  '

# Row 10
$ws.Range("D10").Value = '34.34'
$ws.Range("E10").Value = '  -2.93%  '

# Row 11
$ws.Range("D11").Value = '19.11'
$ws.Range("E11").Value = '  +2.01%  '

# Row 12
$ws.Range("D12").Value = '0.0786'
$ws.Range("E12").Value = '  -0.31%  '

# Row 13
$ws.Range("E13").Value = '  +3.70%  '

# Row 14
$ws.Range("D14").Value = '6.75'
$ws.Range("E14").Value = '  -2.13%  '

# Row 15
$ws.Range("D15").Value = '2.696.82'
$ws.Range("E15").Value = '  +1.16%  '

# Row 16
$ws.Range("D16").Value = '2.332.63'
$ws.Range("E16").Value = '  +0.65%  '

# Row 17
$ws.Range("D17").Value = '0.793'

# Row 18
$ws.Range("D18").Value = '42.967.29'

# Row 19
$ws.Range("D19").Value = '12.20'
$ws.Range("E19").Value = '  -3.40%  '

# Row 20
$ws.Range("D20").Value = '6.19'
$ws.Range("E20").Value = '  +2.37%  '

# Row 21
$ws.Range("E21").Value = '  -0.52%  '

# Row 22
$ws.Range("D22").Value = '67.95'
$ws.Range("E22").Value = '  +0.42%  '

# Row 23
$ws.Range("D23").Value = '237.20'
$ws.Range("E23").Value = '  +0.05%  '

# Row 24
$ws.Range("D24").Value = '2.25'
$ws.Range("E24").Value = '  +3.96%  '

# Row 25
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").Value = '2.43'
$ws.Range("E25").Value = '  +0.31%  '

# Row 26
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.04%  '

# Row 27
$ws.Range("D27").Value = '24.64'
$ws.Range("E27").Value = '  -1.26%  '

# Row 28
$ws.Range("D28").Value = '2.05'
$ws.Range("E28").Value = '  -5.97%  '

# Row 29
$ws.Range("D29").Value = '9.15'
$ws.Range("E29").Value = '  +1.04%  '

# Row 30
$ws.Range("D30").Value = '31.65'
$ws.Range("E30").Value = '  -3.59%  '

# Row 31
$ws.Range("B31").Value = 'FirstDigitalUSD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.10%  '

# Row 32
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").Value = '138.48'
$ws.Range("E32").Value = '  -16.52%  '

# Row 33
$ws.Range("D33").Value = '5.02'
$ws.Range("E33").Value = '  +0.90%  '

# Row 34
$ws.Range("D34").Value = '17.83'
$ws.Range("E34").Value = '  -2.07%  '

# Row 35
$ws.Range("D35").Value = '0.0703'
$ws.Range("E35").Value = '  +1.89%  '

# Row 36
$ws.Range("D36").Value = '4.41'
$ws.Range("E36").Value = '  -1.49%  '

# Row 37
$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").Value = '2.31'
$ws.Range("E37").Value = '  -1.48%  '

# Row 38
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").Value = '1.81'
$ws.Range("E38").Value = '  +2.82%  '

# Row 39
$ws.Range("E39").Value = '  +0.81%  '

# Row 40
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").Value = '2.75'
$ws.Range("E40").Value = '  +0.36%  '

# Row 41
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = '22.31'
$ws.Range("E41").Value = '  +24.18%  '

# Row 42
$ws.Range("E42").Value = '  -0.38%  '

# Row 43
$ws.Range("D43").Value = '1.935.18'
$ws.Range("E43").Value = '  -3.17%  '

# Row 44
$ws.Range("E44").Value = '  -0.11%  '

# Row 45
$ws.Range("D45").Value = '10.07'
$ws.Range("E45").Value = '  -2.69%  '

# Row 46
$ws.Range("E46").Value = '  -2.73%  '

# Row 47
$ws.Range("E47").Value = '  -0.75%  '

# Row 48
$ws.Range("D48").Value = '2.564.41'

# Row 49
$ws.Range("E49").Value = '  +1.47%  '

# Row 50
$ws.Range("D50").Value = '53.63'
$ws.Range("E50").Value = '  +0.31%  '

# Row 51
$ws.Range("D51").Value = '73.13'
$ws.Range("E51").Value = '  +1.72%  '
